$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 473.63635
$ws.Range("I12").Value = 440
$ws.Range("J12").Value = 563.3333
$ws.Range("K12").Value = 440
$ws.Range("L12").Value = 563.3333
$ws.Range("M12").Value = -270
$ws.Range("N12").Value = -903.3333
$ws.Range("H62").Value = 8640.583000000001
$ws.Range("I62").Value = 7376
$ws.Range("K62").Value = 7376
$ws.Range("M62").Value = -6752
$ws.Range("H64").Value = 7654.1816
$ws.Range("J64").Value = 9214.143
$ws.Range("L64").Value = 9214.143
$ws.Range("N64").Value = -9710.143
$ws.Range("H65").Value = 8640.583000000001
$ws.Range("I65").Value = 7376
$ws.Range("K65").Value = 36880
$ws.Range("M65").Value = -33760
$ws.Range("H67").Value = 7654.1816
$ws.Range("J67").Value = 9214.143
$ws.Range("L67").Value = 9214.143
$ws.Range("N67").Value = -10930.143
$ws.Range("H69").Value = 9213.244000000001
$ws.Range("J69").Value = 9208.6
$ws.Range("L69").Value = 27625.8
$ws.Range("N69").Value = -29373.8
$ws.Range("H72").Value = 9213.244000000001
$ws.Range("J72").Value = 9208.6
$ws.Range("L72").Value = 82877.40000000001
$ws.Range("N72").Value = -91613.40000000001
$ws.Range("H76").Value = 4329.125
$ws.Range("I76").Value = 4272.1665
$ws.Range("K76").Value = 4272.1665
$ws.Range("M76").Value = -3957.1665
$ws.Range("H79").Value = 4329.125
$ws.Range("I79").Value = 4272.1665
$ws.Range("K79").Value = 4272.1665
$ws.Range("M79").Value = -3180.1665
$ws.Range("H86").Value = 321572350
$ws.Range("I86").Value = 333334720
$ws.Range("J86").Value = 312750560
$ws.Range("K86").Value = 333334720
$ws.Range("L86").Value = 312750560
$ws.Range("M86").Value = -333333597
$ws.Range("N86").Value = -312752806
$ws.Range("H87").Value = 160934.25
$ws.Range("J87").Value = 159583.33
$ws.Range("L87").Value = 159583.33
$ws.Range("N87").Value = -162079.33
$ws.Range("H89").Value = 321572350
$ws.Range("I89").Value = 333334720
$ws.Range("J89").Value = 312750560
$ws.Range("K89").Value = 1666673600
$ws.Range("L89").Value = 1563752800
$ws.Range("M89").Value = -1666667984
$ws.Range("N89").Value = -1563764032
$ws.Range("H90").Value = 160934.25
$ws.Range("J90").Value = 159583.33
$ws.Range("L90").Value = 478749.99
$ws.Range("N90").Value = -491229.99
$ws.Range("H96").Value = 1901.9231
$ws.Range("I96").Value = 1477.0834
$ws.Range("J96").Value = 7000
$ws.Range("K96").Value = 4431.2502
$ws.Range("L96").Value = 21000
$ws.Range("M96").Value = -3058.2502
$ws.Range("N96").Value = -23746
$ws.Range("H100").Value = 2730
$ws.Range("I100").Value = 1311.125
$ws.Range("J100").Value = 3865.1
$ws.Range("K100").Value = 1311.125
$ws.Range("L100").Value = 3865.1
$ws.Range("M100").Value = -770.125
$ws.Range("N100").Value = -4947.1
$ws.Range("H101").Value = 1502.8
$ws.Range("I101").Value = 1542.3334
$ws.Range("J101").Value = 1443.5
$ws.Range("K101").Value = 4627.0002
$ws.Range("L101").Value = 4330.5
$ws.Range("M101").Value = -3005.0002
$ws.Range("N101").Value = -7574.5
$ws.Range("H103").Value = 870.4583
$ws.Range("I103").Value = 839.6842
$ws.Range("J103").Value = 987.4
$ws.Range("K103").Value = 2519.0526
$ws.Range("L103").Value = 2962.2
$ws.Range("M103").Value = -1933.0526
$ws.Range("N103").Value = -4134.2
$ws.Range("H132").Value = 2891.543
$ws.Range("I132").Value = 3070.459
$ws.Range("K132").Value = 9211.377
$ws.Range("M132").Value = -6681.377
$ws.Range("H135").Value = 916.2432
$ws.Range("I135").Value = 1012.42426
$ws.Range("K135").Value = 9111.81834
$ws.Range("M135").Value = -6576.81834
$ws.Range("H137").Value = 2746.9333
$ws.Range("I137").Value = 2799.6
$ws.Range("K137").Value = 8398.799999999999
$ws.Range("M137").Value = -5848.799999999999
$ws.Range("H138").Value = 3838.6702
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3838.6702
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 11516.0106
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -21796.0106

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 13106815
$ws.Range("I32").Value = 13198902
$ws.Range("K32").Value = 13198902
$ws.Range("M32").Value = -13198615
$ws.Range("H42").Value = 28514
$ws.Range("I42").Value = 30028
$ws.Range("J42").Value = 27000
$ws.Range("K42").Value = 30028
$ws.Range("L42").Value = 27000
$ws.Range("M42").Value = -29542
$ws.Range("N42").Value = -27972
$ws.Range("H45").Value = 4444.727
$ws.Range("I45").Value = 4321.3335
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 4321.3335
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -3944.3335
$ws.Range("N45").Value = -5754
$ws.Range("H61").Value = 1872.0344
$ws.Range("I61").Value = 1791.64
$ws.Range("K61").Value = 1791.64
$ws.Range("M61").Value = -1579.64
$ws.Range("H63").Value = 4116.6665
$ws.Range("I63").Value = 2750
$ws.Range("J63").Value = 4800
$ws.Range("K63").Value = 2750
$ws.Range("L63").Value = 4800
$ws.Range("M63").Value = -2064
$ws.Range("N63").Value = -6172
$ws.Range("H66").Value = 4116.6665
$ws.Range("I66").Value = 2750
$ws.Range("J66").Value = 4800
$ws.Range("K66").Value = 13750
$ws.Range("L66").Value = 24000
$ws.Range("M66").Value = -10318
$ws.Range("N66").Value = -30864
$ws.Range("H74").Value = 1588.5652
$ws.Range("I74").Value = 1545.8422
$ws.Range("K74").Value = 1545.8422
$ws.Range("M74").Value = -671.8422
$ws.Range("H77").Value = 1588.5652
$ws.Range("I77").Value = 1545.8422
$ws.Range("K77").Value = 7729.211
$ws.Range("M77").Value = -3361.211
$ws.Range("H82").Value = 36633.332
$ws.Range("J82").Value = 36633.332
$ws.Range("L82").Value = 36633.332
$ws.Range("N82").Value = -37355.332
$ws.Range("H85").Value = 36633.332
$ws.Range("J85").Value = 36633.332
$ws.Range("L85").Value = 36633.332
$ws.Range("N85").Value = -39129.332
$ws.Range("H97").Value = 1719.8
$ws.Range("I97").Value = 1566.3334
$ws.Range("J97").Value = 1950
$ws.Range("K97").Value = 1566.3334
$ws.Range("L97").Value = 1950
$ws.Range("M97").Value = -1070.3334
$ws.Range("N97").Value = -2942
$ws.Range("H102").Value = 3465.6667
$ws.Range("I102").Value = 1397
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 1397
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = 225
$ws.Range("N102").Value = -7744
$ws.Range("H109").Value = 64084.5
$ws.Range("J109").Value = 64084.5
$ws.Range("L109").Value = 64084.5
$ws.Range("N109").Value = -66858.5
$ws.Range("H125").Value = 86607.414
$ws.Range("J125").Value = 96928.89999999999
$ws.Range("L125").Value = 96928.89999999999
$ws.Range("N125").Value = -106768.9
$ws.Range("H136").Value = 1872.0344
$ws.Range("I136").Value = 1791.64
$ws.Range("K136").Value = 5374.92
$ws.Range("M136").Value = -2824.92

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
$ws.Range("H94").Value = 1151.2195
$ws.Range("I94").Value = 886.4643
$ws.Range("J94").Value = 1721.4615
$ws.Range("K94").Value = 886.4643
$ws.Range("L94").Value = 1721.4615
$ws.Range("M94").Value = -435.4643
$ws.Range("N94").Value = -2623.4615
$ws.Range("H105").Value = 2531.9583
$ws.Range("I105").Value = 2153.2778
$ws.Range("J105").Value = 3668
$ws.Range("K105").Value = 2153.2778
$ws.Range("L105").Value = 3668
$ws.Range("M105").Value = -406.2777999999998
$ws.Range("N105").Value = -7162
$ws.Range("H117").Value = 69920
$ws.Range("J117").Value = 69920
$ws.Range("L117").Value = 69920
$ws.Range("N117").Value = -79098

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2484.0513
$ws.Range("I31").Value = 1864.8695
$ws.Range("J31").Value = 3374.125
$ws.Range("K31").Value = 1864.8695
$ws.Range("L31").Value = 3374.125
$ws.Range("M31").Value = -1569.8695
$ws.Range("N31").Value = -3964.125
$ws.Range("H34").Value = 2484.0513
$ws.Range("I34").Value = 1864.8695
$ws.Range("J34").Value = 3374.125
$ws.Range("K34").Value = 1864.8695
$ws.Range("L34").Value = 3374.125
$ws.Range("M34").Value = -1662.8695
$ws.Range("N34").Value = -3778.125
$ws.Range("H62").Value = 4499.5
$ws.Range("I62").Value = 4499.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4499.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3875.5
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 4499.5
$ws.Range("I65").Value = 4499.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22497.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19377.5
$ws.Range("N65").Value = ""
$ws.Range("H86").Value = 47514.777
$ws.Range("I86").Value = 48688.625
$ws.Range("J86").Value = 47179.395
$ws.Range("K86").Value = 48688.625
$ws.Range("L86").Value = 47179.395
$ws.Range("M86").Value = -47565.625
$ws.Range("N86").Value = -49425.395
$ws.Range("H89").Value = 47514.777
$ws.Range("I89").Value = 48688.625
$ws.Range("J89").Value = 47179.395
$ws.Range("K89").Value = 243443.125
$ws.Range("L89").Value = 235896.975
$ws.Range("M89").Value = -237827.125
$ws.Range("N89").Value = -247128.975
$ws.Range("H99").Value = 1998
$ws.Range("I99").Value = 1998
$ws.Range("K99").Value = 1998
$ws.Range("M99").Value = -500
$ws.Range("H105").Value = 1529.375
$ws.Range("I105").Value = 1247.8572
$ws.Range("K105").Value = 1247.8572
$ws.Range("M105").Value = 499.1428000000001
$ws.Range("H120").Value = 38999.4
$ws.Range("J120").Value = 38749.5
$ws.Range("L120").Value = 38749.5
$ws.Range("N120").Value = -46007.5
$ws.Range("H122").Value = 11116518
$ws.Range("I122").Value = 14291564
$ws.Range("J122").Value = 3857
$ws.Range("K122").Value = 42874692
$ws.Range("L122").Value = 11571
$ws.Range("M122").Value = -42872242
$ws.Range("N122").Value = -16471
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H126").Value = 1998
$ws.Range("I126").Value = 1998
$ws.Range("K126").Value = 5994
$ws.Range("M126").Value = -3524
$ws.Range("H132").Value = 53105.227
$ws.Range("I132").Value = 31542.715
$ws.Range("K132").Value = 94628.145
$ws.Range("M132").Value = -92098.145

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H12").Value = 29.583334
$ws.Range("J12").Value = 27
$ws.Range("L12").Value = 81
$ws.Range("N12").Value = -427
$ws.Range("H120").Value = 23642.785
$ws.Range("J120").Value = 24666.555
$ws.Range("L120").Value = 73999.66500000001
$ws.Range("N120").Value = -83675.66500000001
$ws.Range("H122").Value = 774983.4399999999
$ws.Range("I122").Value = 393.57144
$ws.Range("J122").Value = 3486048
$ws.Range("K122").Value = 3542.14296
$ws.Range("L122").Value = 31374432
$ws.Range("M122").Value = -1092.14296
$ws.Range("N122").Value = -31379332
$ws.Range("H128").Value = 167028
$ws.Range("I128").Value = 167028
$ws.Range("K128").Value = 501084
$ws.Range("M128").Value = -496104
$ws.Range("H131").Value = 1695
$ws.Range("I131").Value = 888.3333
$ws.Range("K131").Value = 2664.9999
$ws.Range("M131").Value = 2375.0001

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 3368.3333
$ws.Range("I80").Value = 3052.5
$ws.Range("K80").Value = 3052.5
$ws.Range("M80").Value = -2054.5
$ws.Range("H83").Value = 3368.3333
$ws.Range("I83").Value = 3052.5
$ws.Range("K83").Value = 15262.5
$ws.Range("M83").Value = -10270.5
$ws.Range("H113").Value = 1236.4445
$ws.Range("I113").Value = 1223.4
$ws.Range("J113").Value = 1252.75
$ws.Range("K113").Value = 1223.4
$ws.Range("L113").Value = 1252.75
$ws.Range("M113").Value = 946.5999999999999
$ws.Range("N113").Value = -5592.75
$ws.Range("H126").Value = 3253.25
$ws.Range("I126").Value = 3171
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 9513
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -7043
$ws.Range("N126").Value = -15440
$ws.Range("H134").Value = 164999.5
$ws.Range("J134").Value = 164999.5
$ws.Range("L134").Value = 494998.5
$ws.Range("N134").Value = -500068.5
$ws.Range("H136").Value = 79299.125
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 79299.125
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 237897.375
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -242997.375

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 9837.166999999999
$ws.Range("I7").Value = 9699.25
$ws.Range("J7").Value = 9947.5
$ws.Range("K7").Value = 9699.25
$ws.Range("L7").Value = 9947.5
$ws.Range("M7").Value = -9587.25
$ws.Range("N7").Value = -10171.5
$ws.Range("H12").Value = 8798.799999999999
$ws.Range("J12").Value = 8798.799999999999
$ws.Range("L12").Value = 8798.799999999999
$ws.Range("N12").Value = -9138.799999999999
$ws.Range("H40").Value = 12350382
$ws.Range("I40").Value = 15876685
$ws.Range("J40").Value = 8324
$ws.Range("K40").Value = 15876685
$ws.Range("L40").Value = 8324
$ws.Range("M40").Value = -15876549
$ws.Range("N40").Value = -8596
$ws.Range("H46").Value = 5135.156
$ws.Range("J46").Value = 3177.1614
$ws.Range("L46").Value = 3177.1614
$ws.Range("N46").Value = -3553.1614
$ws.Range("H55").Value = 379.95456
$ws.Range("I55").Value = 261.18182
$ws.Range("K55").Value = 261.18182
$ws.Range("M55").Value = -88.18182000000002
$ws.Range("H68").Value = 6590.9443
$ws.Range("I68").Value = 7612
$ws.Range("J68").Value = 5569.8887
$ws.Range("K68").Value = 7612
$ws.Range("L68").Value = 5569.8887
$ws.Range("M68").Value = -6863
$ws.Range("N68").Value = -7067.8887
$ws.Range("H71").Value = 6590.9443
$ws.Range("I71").Value = 7612
$ws.Range("J71").Value = 5569.8887
$ws.Range("K71").Value = 38060
$ws.Range("L71").Value = 27849.4435
$ws.Range("M71").Value = -34316
$ws.Range("N71").Value = -35337.4435
$ws.Range("H93").Value = 43479668
$ws.Range("I93").Value = 66667524
$ws.Range("J93").Value = 2436.75
$ws.Range("K93").Value = 66667524
$ws.Range("L93").Value = 2436.75
$ws.Range("M93").Value = -66666276
$ws.Range("N93").Value = -4932.75
$ws.Range("H100").Value = 4250.875
$ws.Range("J100").Value = 5001.4
$ws.Range("L100").Value = 5001.4
$ws.Range("N100").Value = -6083.4
$ws.Range("H126").Value = 9837.166999999999
$ws.Range("I126").Value = 9699.25
$ws.Range("J126").Value = 9947.5
$ws.Range("K126").Value = 29097.75
$ws.Range("L126").Value = 29842.5
$ws.Range("M126").Value = -26627.75
$ws.Range("N126").Value = -34782.5
$ws.Range("H127").Value = 25645.834
$ws.Range("J127").Value = 25645.834
$ws.Range("L127").Value = 25645.834
$ws.Range("N127").Value = -35565.834
$ws.Range("H132").Value = 4572.136
$ws.Range("I132").Value = 4456.778
$ws.Range("K132").Value = 13370.334
$ws.Range("M132").Value = -10840.334
$ws.Range("H136").Value = 2030.0714
$ws.Range("I136").Value = 1674
$ws.Range("J136").Value = 3098.2856
$ws.Range("K136").Value = 5022
$ws.Range("L136").Value = 9294.856800000001
$ws.Range("M136").Value = -2472
$ws.Range("N136").Value = -14394.8568

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5346
$ws.Range("H51").Value = 43530.668
$ws.Range("I51").Value = 27063.334
$ws.Range("K51").Value = 27063.334
$ws.Range("M51").Value = -26553.334
$ws.Range("H52").Value = 40614.2
$ws.Range("I52").Value = 27691.666
$ws.Range("J52").Value = 59998
$ws.Range("K52").Value = 27691.666
$ws.Range("L52").Value = 59998
$ws.Range("M52").Value = -27465.666
$ws.Range("N52").Value = -60450
$ws.Range("H55").Value = 13737.25
$ws.Range("I55").Value = 4999
$ws.Range("J55").Value = 16650
$ws.Range("K55").Value = 4999
$ws.Range("L55").Value = 16650
$ws.Range("M55").Value = -4722
$ws.Range("N55").Value = -17204
$ws.Range("H62").Value = 6150
$ws.Range("H65").Value = 6150
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H81").Value = 4904.9287
$ws.Range("I81").Value = 4733.6
$ws.Range("K81").Value = 9467.200000000001
$ws.Range("M81").Value = -8406.200000000001
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H84").Value = 4904.9287
$ws.Range("I84").Value = 4733.6
$ws.Range("K84").Value = 47336
$ws.Range("M84").Value = -42032
$ws.Range("H126").Value = 3795.2903
$ws.Range("I126").Value = 4026.16
$ws.Range("J126").Value = 2833.3333
$ws.Range("K126").Value = 12078.48
$ws.Range("L126").Value = 8499.999899999999
$ws.Range("M126").Value = -9608.48
$ws.Range("N126").Value = -13439.9999
$ws.Range("H132").Value = 2122.3699
$ws.Range("I132").Value = 1913.6296
$ws.Range("K132").Value = 5740.8888
$ws.Range("M132").Value = -3210.8888
$ws.Range("H136").Value = 25614.166
$ws.Range("I136").Value = 1378.3462
$ws.Range("J136").Value = 64997.375
$ws.Range("K136").Value = 4135.0386
$ws.Range("L136").Value = 194992.125
$ws.Range("M136").Value = -1585.0386
$ws.Range("N136").Value = -200092.125
